# Insert a new data row at row 87 (pushes existing rows 87..197 down to 88..198),
# then populate the new row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(87).Insert()

$ws.Range("A87").Value = 7
$ws.Range("B87").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C87").Value = "Ñuble"
$ws.Range("D87").Value = (Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D87").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E87").Value = 16
$ws.Range("F87").Value = 100114001
$ws.Range("G87").Value = "Papa"
$ws.Range("H87").Value = "Patagonia"
$ws.Range("I87").Value = "1a (guarda)"
$ws.Range("J87").Value = 160
$ws.Range("K87").Value = 7000
$ws.Range("L87").Value = 7500
$ws.Range("M87").Value = 7250
$ws.Range("N87").Value = '$/saco 25 kilos'
$ws.Range("O87").Value = "Provincia de Diguillín"
$ws.Range("P87").Value = 290
$ws.Range("Q87").Value = 25
$ws.Range("R87").Value = "Hortaliza"
